# Aggiornamento fino a 21 marzo
# Appends 4 new daily rows (230-233) with per-comune "nuovi positivi" data
# for dates 2021-04-18 .. 2021-04-21 (Excel serials 44304-44307), mirroring
# the existing layout: column A holds the date (same style as the rows
# above it), columns B..AX hold the per-comune counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data block: date serial, then comune columns B..AX (50 values per row).
$rowsData = @(
    "44304,15,6,5,95,68,26,27,12,15,1,23,48,48,2,4,3,30,5,12,43,260,9,12,11,35,4,9,14,7,33,7,126,8,23,3,27,23,52,4,39,1221,19,2,3,1,1,2,1,0",
    "44305,15,4,6,83,67,17,18,11,13,1,21,58,48,2,3,1,34,7,12,43,253,9,7,11,33,1,8,13,7,41,6,130,9,23,2,27,20,49,2,40,1169,5,0,3,2,1,2,1,0",
    "44306,15,5,2,79,70,16,22,10,12,1,16,51,44,2,2,0,29,7,12,39,260,9,8,16,40,1,8,15,7,38,5,127,11,31,4,26,27,51,2,21,1154,3,0,3,2,1,2,2,0",
    "44307,21,5,4,78,73,15,22,10,13,1,15,52,40,2,1,0,30,7,12,38,259,9,6,17,39,3,12,14,8,37,5,127,10,30,5,27,27,52,2,24,1165,3,0,3,1,1,3,2,0"
)

$startRow = 230

for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $targetRow = $startRow + $i
    $values = $rowsData[$i].Split(",")

    # Column A: carry the date-cell style forward from the row above
    # (format-only paste), then write this row's date serial value.
    $ws.Range("A" + ($targetRow - 1)).Copy()
    $ws.Range("A" + $targetRow).PasteSpecial(-4122)
    $ws.Cells.Item($targetRow, 1).Value = [double]$values[0]

    # Columns B..AX: plain numeric counts, one cell at a time.
    for ($col = 2; $col -le $values.Count; $col++) {
        $ws.Cells.Item($targetRow, $col).Value = [double]$values[$col - 1]
    }
}

$excel.CutCopyMode = $false

Write-Host "Added rows 230-233"
